$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new vocabulary -----------------------------------------------
# Three of the new word pairs were typed in before the list got re-sorted,
# so they belong inside the alphabetically sorted block. Enter them first
# (order matters: it drives the order brand-new strings are appended to
# the shared-strings table).
$ws.Range("A35").Value = "apartment"
$ws.Range("B35").Value = "căn hộ"
$ws.Range("A36").Value = "document"
$ws.Range("B36").Value = "tài liệu"
$ws.Range("A37").Value = "electricity"
$ws.Range("B37").Value = "điện"

# Fix the header capitalisation (Eng -> ENG). This must happen after the
# three pairs above and before the remaining new vocabulary below.
$ws.Range("A1").Value = "ENG"

# --- Sort the table alphabetically by the English word -----------------
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A1"))
$sort.SetRange($ws.Range("A1:B37"))
$sort.Header = 1
$sort.Apply()

# --- Append the rest of the new vocabulary after the sorted block ------
# These words were added after the list had already been sorted, so they
# stay in entry order at the bottom of the sheet.
$newWords = @(
    @("jar", "cái lọ"),
    @("purse", "cái ví"),
    @("strength", "sức mạnh"),
    @("vegetable", "rau củ"),
    @("mathematician", "nhà toán học"),
    @("physicist", "nhà vật lý"),
    @("scientist", "nhà khoa học"),
    @("circuit", "bo mạch"),
    @("napkin", "khăn ăn"),
    @("construction", "công trường"),
    @("saucepan", "cái chảo"),
    @("cupboard", "tủ đựng chén"),
    @("chopsticks", "đôi đũa"),
    @("spoon", "cái muỗng")
)

$row = 38
foreach ($pair in $newWords) {
    $ws.Range("A$row").Value = $pair[0]
    $ws.Range("B$row").Value = $pair[1]
    $row = $row + 1
}

# --- Cosmetic touch-ups to match the refreshed look of the sheet -------
# Switch the whole workbook's base font from Calibri to Arial.
$wb.Styles.Item("Normal").Font.Name = "Arial"

# Resize the two columns to fit their (now wider) contents.
$ws.Range("A:B").AutoFit() | Out-Null

# Leave the selection on the last populated cell, like the author did.
$ws.Range("A51").Select() | Out-Null
